# DevTesting_IC.dry.xlsx edits
# - Edit IC to include release tier inputs (CoordinatedOps!C3:C5)
# - Make CoordinatedOps the active/selected sheet & update its selection
# - Update Reservoirs sheet view (no longer the selected tab; pane scroll +
#   selection change)
# - Update InterveningFlow sheet selection

$wb = $excel.ActiveWorkbook

# --- CoordinatedOps: release tier inputs (IC) -----------------------------
$wsCO = $wb.Worksheets.Item("CoordinatedOps")
$wsCO.Range("C3").Value = 3
$wsCO.Range("C4").Value = 3
$wsCO.Range("C5").Value = 3

# --- Reservoirs: scroll the frozen pane & move the selection --------------
$wsRes = $wb.Worksheets.Item("Reservoirs")
$wsRes.Activate()
$wsRes.Range("A17:XFD17").Select()

# --- InterveningFlow: move the selection -----------------------------------
$wsIF = $wb.Worksheets.Item("InterveningFlow")
$wsIF.Activate()
$wsIF.Range("B33:I33").Select()

# --- CoordinatedOps becomes the active tab with C6 selected ----------------
$wsCO.Activate()
$wsCO.Range("C6").Select()
